$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at the top of the Chirimoya weekly block (row 263) to make room
# for the newest week (2023-09-25, serial 45194), pushing all subsequent weekly
# rows down by 3 (263-299 -> 266-302).
$ws.Rows("263:265").Insert()

# Populate the new week (Especial / Primera / Segunda) with its price data.
$ws.Cells.Item(263,1).Value = 8
$ws.Cells.Item(263,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(263,3).Value = "Coquimbo"
$ws.Cells.Item(263,4).Value = 45194
$ws.Cells.Item(263,5).Value = 4
$ws.Cells.Item(263,6).Value = "Fruta"
$ws.Cells.Item(263,7).Value = 100107
$ws.Cells.Item(263,8).Value = "Otros"
$ws.Cells.Item(263,9).Value = 100107002
$ws.Cells.Item(263,10).Value = "Chirimoya"
$ws.Cells.Item(263,11).Value = "Cultivar IV Región"
$ws.Cells.Item(263,12).Value = "Especial"
$ws.Cells.Item(263,13).Value = 240
$ws.Cells.Item(263,14).Value = 23000
$ws.Cells.Item(263,15).Value = 24000
$ws.Cells.Item(263,16).Value = 23500
$ws.Cells.Item(263,17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(263,18).Value = "Provincia de Limarí"
$ws.Cells.Item(263,19).Value = 2350
$ws.Cells.Item(263,20).Value = 10
$ws.Cells.Item(264,1).Value = 8
$ws.Cells.Item(264,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(264,3).Value = "Coquimbo"
$ws.Cells.Item(264,4).Value = 45194
$ws.Cells.Item(264,5).Value = 4
$ws.Cells.Item(264,6).Value = "Fruta"
$ws.Cells.Item(264,7).Value = 100107
$ws.Cells.Item(264,8).Value = "Otros"
$ws.Cells.Item(264,9).Value = 100107002
$ws.Cells.Item(264,10).Value = "Chirimoya"
$ws.Cells.Item(264,11).Value = "Cultivar IV Región"
$ws.Cells.Item(264,12).Value = "Primera"
$ws.Cells.Item(264,13).Value = 300
$ws.Cells.Item(264,14).Value = 21000
$ws.Cells.Item(264,15).Value = 22000
$ws.Cells.Item(264,16).Value = 21500
$ws.Cells.Item(264,17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(264,18).Value = "Provincia de Limarí"
$ws.Cells.Item(264,19).Value = 2150
$ws.Cells.Item(264,20).Value = 10
$ws.Cells.Item(265,1).Value = 8
$ws.Cells.Item(265,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(265,3).Value = "Coquimbo"
$ws.Cells.Item(265,4).Value = 45194
$ws.Cells.Item(265,5).Value = 4
$ws.Cells.Item(265,6).Value = "Fruta"
$ws.Cells.Item(265,7).Value = 100107
$ws.Cells.Item(265,8).Value = "Otros"
$ws.Cells.Item(265,9).Value = 100107002
$ws.Cells.Item(265,10).Value = "Chirimoya"
$ws.Cells.Item(265,11).Value = "Cultivar IV Región"
$ws.Cells.Item(265,12).Value = "Segunda"
$ws.Cells.Item(265,13).Value = 200
$ws.Cells.Item(265,14).Value = 18000
$ws.Cells.Item(265,15).Value = 19000
$ws.Cells.Item(265,16).Value = 18500
$ws.Cells.Item(265,17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(265,18).Value = "Provincia de Limarí"
$ws.Cells.Item(265,19).Value = 1850
$ws.Cells.Item(265,20).Value = 10
